# Update the "Login" worksheet with two additional rows of test data,
# then move the active selection to C4 (matching the saved workbook view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$ws.Range("A2").Value = "skv"
$ws.Range("B2").Value = "shubham"
$ws.Range("A3").Value = "shubham"
$ws.Range("B3").Value = "shubham"

$ws.Range("C4").Select()
